# Commit 2019-08-09 kl. 09:20
# Update the "Fixzone" tracking sheet: one more day logged for week ending
# row 21 (I21: 3 -> 4), which flows through the K3/L3/M3 summary formulas
# automatically on recalculation. Also fill in cell G21 to match the rest
# of the row (green "done" highlight) and leave the selection where the
# user was last working (L20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the day count for this row from 3 to 4.
$ws.Range("I21").Value = 4

# G21 was missing the green fill that the rest of the row (D21:F21)
# already has; match it so the row is visually consistent.
$ws.Range("G21").Interior.Color = $ws.Range("D21").Interior.Color

# Leave the selection on L20, matching where the edit was made.
$ws.Range("L20").Select() | Out-Null
